# Add a new "section_title" column (B) that shows "Стол <section_id>" for
# each question row, and format the section_id column (A) as an integer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember column A's width before we insert, so the new column B can match it.
$sectionIdColWidth = $ws.Columns("A").ColumnWidth

# Insert a new, blank column before the old column B (question text), shifting
# every column from B onward one slot to the right.
$ws.Columns("B:B").Insert()

# Column A (section_id) should display as a plain integer.
$ws.Range("A2:A24").NumberFormat = "0"

# New column B mirrors column A's width and gets the header + derived values.
$ws.Columns("B").ColumnWidth = $sectionIdColWidth

$ws.Range("B1").Value = "section_title"

for ($r = 2; $r -le 24; $r++) {
    $sectionId = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = "Стол " + $sectionId
}

# Match the view state captured in the saved workbook: zoomed out a bit, with
# the last section-title cell selected.
$excel.ActiveWindow.Zoom = 70
$ws.Range("B24").Select()
